$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three "Ista (sub-sample)" cells to "Ista (EPC sample)"
$ws.Range("E12").Value = "Ista (EPC sample)"
$ws.Range("E13").Value = "Ista (EPC sample)"
$ws.Range("E14").Value = "Ista (EPC sample)"

# Update the saved selection to match the diff (E12:E14, active cell E12)
$ws.Range("E12:E14").Select()
